# A new data row (for a later "Puerro" price observation) was inserted
# into the daily log at row 131. All rows that previously occupied
# 131-220 shift down by one (132-221); the worksheet dimension grows
# from A1:R220 to A1:R221 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 131, pushing existing rows 131:220
# down to 132:221. Excel copies formatting from the row above, which
# gives the new D131 cell the same date style (s="2") as the rest of
# column D.
$ws.Rows("131:131").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A131").Value = 10
$ws.Range("B131").Value = "Vega Modelo de Temuco"
$ws.Range("C131").Value = "La Araucanía"
$ws.Range("D131").Value = 44767
$ws.Range("E131").Value = 9
$ws.Range("F131").Value = 100112005
$ws.Range("G131").Value = "Puerro"
$ws.Range("H131").Value = "Azul de Maquehue"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 20
$ws.Range("K131").Value = 16000
$ws.Range("L131").Value = 16000
$ws.Range("M131").Value = 16000
$ws.Range("N131").Value = "`$/docena de paquetes"
$ws.Range("O131").Value = "Provincia de Cautín"
$ws.Range("P131").Value = 1333
$ws.Range("Q131").Value = 12
$ws.Range("R131").Value = "Hortaliza"
